$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update age-bracket labels in every header: "(6-17)" -> "(5-17)"
#    (Children/Girls/Boys, for overall + Host/IDP/Returnees/Refugees breakdowns)
$ws.Cells.Replace("(6-17)", "(5-17)")

# 2. Add the "Host -- Children" status column (F) = 70% of total (C column)
$ws.Range("F2").Formula = "=C2*0.7"
$ws.Range("F3:F25").Formula = "=C3*0.7"

# 3. Add the "IDP -- Girls" status column (J) = 30% of total (C column)
$ws.Range("J2").Formula = "=C2*0.3"
$ws.Range("J3:J25").Formula = "=C3*0.3"

# 4. Update the view state: scroll down a bit and move the active selection
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("M18").Select()
